$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 235
$ws1.Range("F4").Value = 859
$ws1.Range("F6").Value = 35

# Sheet "全部类型" (sheet4): update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 235
$ws4.Range("F5").Value = 859
$ws4.Range("F7").Value = 35
